# Status.xlsx update: add a new status row (18) and refresh the JSF research
# note in row 17 now that admin-note work has landed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 18 - same date style/format as the row above it.
$ws.Range("A17").Copy()
$ws.Range("A18").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A18").Value = (Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0).AddDays(43154)

$ws.Range("B18").Value = "Android display of all entries (100%). Communication of start- and enddate between two activities (30%)"
$ws.Range("C18").Value = "C# LoadAdminGuiElements in combination with webservice (100%). Implementation of adminNote in all windows (40%)"

# Row 17's JSF research note now also covers the follow-up work.
$ws.Range("D17").Value = "Research Web App Frameworks, decision: JSF, reading up in JSF (10%), WebApp Login (100%), Display of all entries (100%), AddEntry (100%)"

$ws.Range("D18").Value = "Reading up in JSF (100%), Restore of MongoDB (it was empty) (100%), MongoDB adminNote added (100%) "

$ws.Range("D18").Select()
